$wb = $excel.ActiveWorkbook

# ---- Sheet: SortedActivities ----
$ws1 = $wb.Worksheets.Item("SortedActivities")

# Insert a new row at position 134 for the newly added "Harishchandragad Trek by e2e" activity
$ws1.Rows.Item(134).Insert()
$ws1.Range("A134").Value = '2025-07-19T21:00:00'
$ws1.Range("B134").Value = '2025-08-31T21:00:00'
$ws1.Range("C134").Value = 'Harishchandragad Trek by e2e'
$ws1.Range("D134").Value = 'Harishchandra fort, Pune'
$ws1.Range("E134").Value = 1499

# Corrections / reorder fixes for other rows in SortedActivities
$ws1.Range("A9").Value = '2025-07-17T18:00:00'
$ws1.Range("B9").Value = '2025-08-17T18:00:00'
$ws1.Range("A18").Value = '2025-07-17T18:00:00'
$ws1.Range("B18").Value = '2025-08-17T18:00:00'
$ws1.Range("C44").Value = 'Hand Built Pottery Date - For Couples and Friends'
$ws1.Range("C45").Value = 'Coil Pottery Date for Couples and Friends'
$ws1.Range("C69").Value = 'Tissue Texture Art Date'
$ws1.Range("A70").Value = '2025-07-19T14:00:00'
$ws1.Range("C70").Value = 'Acrylic Painting Date at Blue Tokai'
$ws1.Range("A71").Value = '2025-07-20T14:00:00'
$ws1.Range("C71").Value = 'Pottery Painting Date for Couples and Friends'
$ws1.Range("C72").Value = 'Texture Art Date | Blue Tokai Coffee Roasters'
$ws1.Range("C73").Value = 'Resin Galaxy Art Date for Couples and Friends | Pune'
$ws1.Range("C74").Value = 'Tissue Texture Art Date - PizzaExpress Pune'
$ws1.Range("C75").Value = 'Resin Art Date for Friends and Couples at PizzaExpress Pune'
$ws1.Range("C76").Value = 'Acrylic Painting Date for Couples and Friends at Pizza Express Pune'
$ws1.Range("C77").Value = 'Resin Glow Art Date at Pizza Express | KOPA Mall'
$ws1.Range("C78").Value = 'Resin Geode Art Date | Pune'
$ws1.Range("C79").Value = 'Paint Your Tote Bag at PizzaExpress Pune'
$ws1.Range("C80").Value = 'Clay Trinket Tray Date - Pizza Express KOPA Mall Pune'
$ws1.Range("C81").Value = 'Clay Miniature Magnet Date- PizzaExpress Pune'
$ws1.Range("C82").Value = 'Fluid Art Date for Couples and Friends at PizzaExpress Pune'
$ws1.Range("C83").Value = 'Resin Koi Pond Art for Couples and Friends | Pune'
$ws1.Range("C84").Value = 'Knife Painting Art for Couples and Friends at Pizza Express Pune'
$ws1.Range("C85").Value = 'Neon Painting Date - PizzaExpress Pune'
$ws1.Range("C86").Value = 'DIY Photo Frame at PizzaExpress Pune'
$ws1.Range("C87").Value = 'Pottery Painting Date - PizzaExpres Pune'
$ws1.Range("C88").Value = 'Resin Beach Art Date- PizzaExpress Pune'
$ws1.Range("C89").Value = 'Canvas Painting Date - PizzaExpress Pune'
$ws1.Range("C90").Value = 'Texture Art Date'
$ws1.Range("B91").Value = '2025-07-20T14:00:00'
$ws1.Range("C91").Value = 'Neon Painting Date'
$ws1.Range("C92").Value = 'Resin Beach Art Date'
$ws1.Range("C94").Value = 'Tissue Texture Art Date'
$ws1.Range("C95").Value = 'Texture Art Date: Create, Connect and Unwind'
$ws1.Range("C96").Value = 'Clay Miniature Magnet Date'
$ws1.Range("B97").Value = '2025-07-27T14:00:00'
$ws1.Range("C97").Value = 'Acrylic Painting Date at Flow Baner'
$ws1.Range("A114").Value = '2025-07-19T05:00:00'
$ws1.Range("B114").Value = '2025-09-21T05:00:00'
$ws1.Range("C114").Value = 'Kalu Waterfall Trek - Trekfit adventures'
$ws1.Range("A115").Value = '2025-07-18T05:00:00'
$ws1.Range("B115").Value = '2025-08-17T05:00:00'
$ws1.Range("C115").Value = 'Devkund Waterfall Trek-Trekfit adventures'
$ws1.Range("A116").Value = '2025-07-19T04:00:00'
$ws1.Range("B116").Value = '2025-09-21T04:00:00'
$ws1.Range("C116").Value = 'Nanemachi Waterfall Trek - Trekfit adventures'
$ws1.Range("C119").Value = 'Labubu Pop Art Date for Couples and Friends'
$ws1.Range("C121").Value = 'Labubu Clay Figurines Workshop for Couples and Friends'
$ws1.Range("C122").Value = 'Japanese Kintsugi Workshop | Pune'
$ws1.Range("C123").Value = 'Labubu Canvas Painting Date for Couples and Friends | Pune'
$ws1.Range("A127").Value = '2025-07-18T11:00:00'
$ws1.Range("B127").Value = '2025-07-30T11:00:00'
$ws1.Range("C143").Value = 'Resin Trinket Tray Date for Couples and Friends'
$ws1.Range("C144").Value = 'Floral Resin Art Date for Couples & Friends - PizzaExpress Pune'
$ws1.Range("C150").Value = 'Waterfall Hike'
$ws1.Range("C154").Value = 'Kids Art Party - PizzaExpress | Pune'
$ws1.Range("C155").Value = 'Couple Art Workshop - PizzaExpress | Pune'
$ws1.Range("C156").Value = 'Date Night Art at PizzaExpress Pune'

# ---- Sheet: SortedEvents ----
$ws2 = $wb.Worksheets.Item("SortedEvents")
$ws2.Range("A6").Value = '20 Jul, 8PM'
$ws2.Range("B6").Value = 'Live Bollywood Music'
$ws2.Range("A7").Value = '19 Jul - 26 Jul, 8PM'
$ws2.Range("B7").Value = 'DJ Night ft. DJ Jack'
$ws2.Range("A26").Value = '20 Jul, 8PM'
$ws2.Range("B26").Value = 'Sunday Night Live by FEELz Comedy'
$ws2.Range("A27").Value = '17 Jul - 18 Jul, 9PM'
$ws2.Range("B27").Value = 'Late Shift Laugh''s Live at AUNDH by FEELz COMEDY'
$ws2.Range("B53").Value = 'Entrepreneurs Meetup by We Founders Collab | Pune'
$ws2.Range("B54").Value = 'Global Startups Club - Startup Networking | Pune'
$ws2.Range("B55").Value = 'Business Networking | Pune'
$ws2.Range("A61").Value = '18 Jul, 10:30PM'
$ws2.Range("B61").Value = 'Friday Night Bounce feat. Life Enjoyers Club'
$ws2.Range("A62").Value = '19 Jul, 10:30PM'
$ws2.Range("B62").Value = 'Sneaky Link feat. Karonik++'

# ---- Sheet: Movies ----
$ws3 = $wb.Worksheets.Item("Movies")
$ws3.Range("A10").Value = 'Aankhon Ki Gustaakhiyan'
$ws3.Range("B10").Value = 'Hindi'
$ws3.Range("A11").Value = 'Jarann'
$ws3.Range("B11").Value = 'Marathi'
$ws3.Range("A13").Value = 'Bhaag Milkha Bhaag (2013)'
$ws3.Range("B13").Value = 'Hindi'
$ws3.Range("C13").Value = 'U'
$ws3.Range("A14").Value = 'How to Train Your Dragon'
$ws3.Range("B14").Value = 'English'
$ws3.Range("C14").Value = 'UA7+'
$ws3.Range("A15").Value = 'Janaki V vs State Of Kerala'
$ws3.Range("B15").Value = 'Malayalam'
$ws3.Range("C15").Value = 'UA16+'
$ws3.Range("A16").Value = 'Ye Re Ye Re Paisa 3'
$ws3.Range("B16").Value = 'Marathi'
$ws3.Range("C16").Value = 'UA13+'
$ws3.Range("A17").Value = 'Ekka'
$ws3.Range("B17").Value = 'Kannada'
$ws3.Range("C17").Value = 'UA16+'
$ws3.Range("A18").Value = 'I Know What You Did Last Summer'
$ws3.Range("B18").Value = 'English'
$ws3.Range("C18").Value = 'A'
$ws3.Range("A19").Value = '5th September'
$ws3.Range("B19").Value = 'Hindi'
$ws3.Range("C19").Value = 'U'
$ws3.Range("A20").Value = 'Junior'
$ws3.Range("B20").Value = 'Telugu'
$ws3.Range("C20").Value = 'UA13+'
$ws3.Range("A21").Value = 'Tanvi The Great'
$ws3.Range("B21").Value = 'Hindi'
$ws3.Range("C21").Value = 'U'
$ws3.Range("A22").Value = 'Nilgiris: A Shared Wilderness'
$ws3.Range("A23").Value = 'Smurfs'
$ws3.Range("A26").Value = 'Nintama Rantaro: Invincible Master Of The Dokutake Ninja'
$ws3.Range("B26").Value = 'Japanese'
$ws3.Range("A27").Value = 'Raas'
$ws3.Range("B27").Value = 'Bengali'
$ws3.Range("C27").Value = 'UA7+'
$ws3.Range("A28").Value = 'Sarbala Ji'
$ws3.Range("B28").Value = 'Punjabi'
$ws3.Range("C28").Value = 'UA13+'
